$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at position 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows(13).Insert()

# 2) The inserted row auto-creates an empty A13 cell (inheriting column A's style); remove it entirely
$ws.Range("A13").Clear()

# 3) Populate B13/C13 with the 'Docentes responsaveis' content, copying number/format from B10/C10
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B13").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("C13").Value = "8554681 - Pedro Felipe Arce Castillo"

# 4) Objetivos: row 10 B/C gets the new long description (replacing the old Pedro Felipe text)
$objText = @'
Aplicar os conceitos fundamentais relacionados aos processos físicos químicos, ampliando o conhecimento termodinâmico dos sistemas, isto é, a definição dos critérios de equilíbrio e de espontaneidade para misturas e reações químicas.
'@
$ws.Range("B10").Value = $objText
$ws.Range("C10").Value = $objText

# 5) Programa resumido: row 14 (old 13) B/C gets new short-syllabus summary
$progResumido = @'
Termodinâmica de soluções. Equilíbrio líquido  vapor. Equilíbrio de fases. Equilíbrio em reações químicas  Equilíbrio químico
'@
$ws.Range("B14").Value = $progResumido
$ws.Range("C14").Value = $progResumido

# 6) Programa: row 16 (old 15) B/C gets the full syllabus text (replacing stray date value)
$programa = @'
1- Termodinâmica de soluções 
1.1- Relações fundamentais entre propriedades 
1.2- O potencial químico 
1.3- Fugacidade e coeficiente de fugacidade 
1.4- A solução Ideal 
1.5- Modelos para a energia de Gibbs 
1.6- Propriedades de mistura 
1.7- Efeitos térmicos em processos de mistura 
2- Equilíbrio liquido  vapor 
2.1- A natureza em equilíbrio 
2.2- A regra das fases. Teorema de Duhem 
2.3- Calculo dos pontos de orvalho e de bolha 
2.4- Calculo de Flash 
3- Equilíbrio de fases 
3.1- Equilíbrio e estabilidade 
3.2- Equilíbrio líquido-líquido 
3.3- Equilíbrio líquido-líquido-vapor 
3.4- Equilíbrio sólido-líquido 
3.5- Equilíbrio sólido-vapor 
3.6- Equilíbrio na adsorção de gases em sólidos 
4- Equilíbrio em reações químicas  Equilíbrio químico 
4.1- A variação de energia de Gibbs padrão e a constante de equilíbrio 
4.2- Efeito da temperatura sobre a constante de equilíbrio 
4.3- Avaliação das constantes de equilíbrio 
4.4- Relação entre as constantes de equilíbrio e a composição 
4.5- Conversões de equilíbrio em reações isoladas
'@
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# 7) Metodo: row 19 (old 18) B/C gets the evaluation method text
$metodo = @'
A avaliação será feita por meio de duas provas escritas (P1 e P2).
'@
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# 8) Criterio: row 20 (old 19) B/C gets the final-grade formula text
$criterio = @'
A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2
'@
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# 9) Norma de recuperacao: row 21 (old 20) B/C gets the recovery-exam text
$recup = @'
A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2
'@
$ws.Range("B21").Value = $recup
$ws.Range("C21").Value = $recup

# 10) Bibliografia: row 22 (old 21) B/C gets the full bibliography text
$biblio = @'
KORETSKY, M. D. Termodinâmica para Engenharia Química. 1 ed. LTC Editora, 2007. 
MORAN, M. J.; SHAPIRO, H. N. Princípios de Termodinâmica para Engenharia. 1 ed. LTC Editora, 2009. 
SANDLER, S. I., Chemical and Engineering Thermodynamics, 3rd ed., John Wiley & Sons, 1999 
SMITH, J.M.; VAN NESS, H.C.; Abott, M. M. Introdução à Termodinâmica da Engenharia Química. 7ª ed. LTC editora, 2007. 
TERRON, L. R. Termodinâmica Química Aplicada. 1 ed. Editora Manole Ltda, 2009. 
VAN WILEN, J. Sonntag, Richard. E. Fundamentos da Termodinâmica Clássica. 6 ed. 2004
'@
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

